# Updates cryptos list values (price/volume columns, and two row swaps)
# to match the latest scrape, per commit "Updated cryptos list on Fri Nov  8 17:30:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "75.946.31"
$ws.Range("E2").Value = "  -0.42%  "
# Row 3
$ws.Range("D3").Value = "2.894.85"
$ws.Range("E3").Value = "  +1.05%  "
# Row 4
$ws.Range("E4").Value = "  +0.06%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.23"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.61%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "588.88"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.68%  "
# Row 7
$ws.Range("E7").Value = "  +0.04%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.543"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.68%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.191"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.81%  "
# Row 10
$ws.Range("D10").Value = "2.896.01"
$ws.Range("E10").Value = "  +1.14%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +11.76%  "
# Row 12
$ws.Range("E12").Value = "  +0.28%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.84"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.33%  "
# Row 14
$ws.Range("D14").Value = "3.431.35"
$ws.Range("E14").Value = "  +1.38%  "
# Row 15
$ws.Range("D15").Value = "75.837.97"
$ws.Range("E15").Value = "  -0.36%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.58"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.25%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000184"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.47%  "
# Row 18
$ws.Range("D18").Value = "2.916.12"
$ws.Range("E18").Value = "  +1.79%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +4.31%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.55"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -6.18%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.18"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -4.67%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.63%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.21"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -5.43%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.53"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.63%  "
# Row 25
$ws.Range("E25").Value = "  -0.06%  "
# Row 26
$ws.Range("D26").Value = "3.045.75"
$ws.Range("E26").Value = "  +1.24%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.17"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.22%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.45"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.47%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000103"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.89%  "
# Row 30
$ws.Range("E30").Value = "  +0.01%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +2.32%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.34"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.44%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "489.01"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -5.63%  "
# Row 34
$ws.Range("E34").Value = "  -1.20%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.01%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.62"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.78%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.85"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.03%  "
# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.382"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +10.58%  "
# Row 39
$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.107"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +19.70%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.69"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.15%  "
# Row 41
$ws.Range("E41").Value = "  +0.01%  "
# Row 42
$ws.Range("E42").Value = "  -9.37%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "176.83"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.09%  "
# Row 44
$ws.Range("E44").Value = "  -5.65%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.61"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.51%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.09"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.78%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.15"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -6.45%  "
# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.570"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.46%  "
# Row 49
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.78"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.41%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -7.13%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.16"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +3.88%  "
